$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of column number -> new value (applies identically to rows 2 and 3)
$colValues = @{
    7 = -1.489514066496164   # column G, was -7.938596491228071
    8 = -1.572890025575448   # column H, was -7.938596491228071
    9 = -1.103580562659846   # column I, was -10.08771929824561
    10 = -1.103580562659846   # column J, was -10.08771929824561
    11 = -10.3   # column K, was -13.3
    12 = -1.317135549872123   # column L, was -11.66666666666667
    21 = 2.04   # column U, was 3.27
    22 = 0.009649952696310311   # column V, was 0.01744
    23 = -0.6602564102564104   # column W, was -0.4586206896551724
    24 = 0.09087725862885132   # column X, was 0.1486669971562242
    25 = -0.7511336688852617   # column Y, was -0.6072876868113966
    26 = 0.1578841106400161   # column Z, was 0.01790482173708183
    27 = -0.1742378356551585   # column AA, was -0.1806188157688079
    28 = 0.08321196620529456   # column AB, was 0.1305359632206022
    29 = -0.2574498018604531   # column AC, was -0.3111547789894101
    30 = 33.3   # column AD, was 37.2
    32 = 33.3   # column AF, was 37.2
    33 = 31.26   # column AG, was 33.93
    34 = 0.1360850020433183   # column AH, was 0.1655540720961282
    35 = 0.6330798479087453   # column AI, was 0.7045454545454546
    36 = 0.1288222203906701   # column AJ, was 0.1532312694756808
    37 = 0.6182753164556961   # column AK, was 0.6850393700787402
    38 = 1.85   # column AL, was 2.54
    39 = 1.375   # column AM, was 1.52
    40 = -4.710042432814709   # column AN, was -3.412844036697248
    41 = -4.664864864864865   # column AO, was -4.52755905511811
    42 = -4.421499292786421   # column AP, was -3.112844036697247
    43 = -6.276363636363637   # column AQ, was -7.565789473684211
}

foreach ($row in 2..3) {
    foreach ($colNum in $colValues.Keys) {
        $ws.Cells.Item($row, $colNum).Value = $colValues[$colNum]
    }
}

Write-Host "Applied capital structure updates to rows 2-3"